$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numeric-looking strings (e.g. "26.024.14",
# "1.00") in the source data. Force the column to Text format first so
# assigning these values does not get auto-coerced into numbers (which
# would silently normalize values like "1.00" -> 1 or "216.46" -> 216.46
# as a Double, breaking the literal text representation).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.024.14'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '1.643.20'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.68%  '
$ws.Range("D5").Value = '216.46'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("E7").Value = '  +0.59%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.0639'
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D13").Value = '1.871.39'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '1.619.26'
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("D16").Value = '0.0₃0767'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").Value = '63.09'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '25.966.06'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").Value = '193.25'
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("D22").Value = '9.94'
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '0.133'
$ws.Range("E24").Value = '  +6.69%  '
$ws.Range("E25").Value = '  +1.11%  '
$ws.Range("D26").Value = '144.71'
$ws.Range("E26").Value = '  +1.50%  '
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").Value = '  -1.59%  '
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("E35").Value = '  +2.36%  '
$ws.Range("D36").Value = '0.904'
$ws.Range("E36").Value = '  -0.68%  '
$ws.Range("D37").Value = '1.135.06'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("D42").Value = '99.76'
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").Value = '1.780.30'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("E45").Value = '  +2.83%  '
$ws.Range("D46").Value = '56.81'
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("E47").Value = '  +2.68%  '
$ws.Range("D48").Value = '1.46'
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("D49").Value = '7.78'
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("D51").Value = '0.0962'
$ws.Range("E51").Value = '  -0.19%  '
